# Weekly price-report update: insert a new "Cilantro" record as row 255
# (Feria Lagunitas de Puerto Montt, fecha 2022-07-12 / serial 44754),
# pushing the existing rows 255-315 down to 256-316.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("255:255").Insert()

$ws.Range("A255").Value = 4
$ws.Range("B255").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C255").Value = "Los Lagos"
$ws.Range("D255").Value = 44754
$ws.Range("E255").Value = 10
$ws.Range("F255").Value = 100112040
$ws.Range("G255").Value = "Cilantro"
$ws.Range("H255").Value = "Sin especificar"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 250
$ws.Range("K255").Value = 14000
$ws.Range("L255").Value = 14000
$ws.Range("M255").Value = 14000
$ws.Range("N255").Value = "`$/caja 36 atados"
$ws.Range("O255").Value = "Región Metropolitana"
$ws.Range("P255").Value = 389
$ws.Range("Q255").Value = 36
$ws.Range("R255").Value = "Hortaliza"
